$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.909.25'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.637.92'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.60%  '
$ws.Range("D5").Value = '''215.15'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '''0.520'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("D8").Value = '''29.01'
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").Value = '''0.0610'
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '1.870.50'
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.588'
$ws.Range("E13").Value = '  +3.99%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.564.22'
$ws.Range("E14").Value = '  -3.32%  '
$ws.Range("D15").Value = '''9.44'
$ws.Range("E15").Value = '  +6.95%  '
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '29.920.17'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '''64.65'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").Value = '''240.84'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '0.0₃0703'
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = '''9.93'
$ws.Range("E22").Value = '  +3.70%  '
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("D24").Value = '''2.18'
$ws.Range("E24").Value = '  +3.06%  '
$ws.Range("D25").Value = '''157.46'
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").Value = '''15.52'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").Value = '''6.63'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("D30").Value = '''0.0490'
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = '''3.39'
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("D33").Value = '''3.20'
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").Value = '1.424.72'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("E35").Value = '  +3.47%  '
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").Value = '''76.45'
$ws.Range("E40").Value = '  +10.60%  '
$ws.Range("D41").Value = '''0.559'
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").Value = '''0.834'
$ws.Range("E42").Value = '  +1.07%  '
$ws.Range("D43").Value = '''0.0499'
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''5.37'
$ws.Range("E47").Value = '  -1.14%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.778.23'
$ws.Range("E48").Value = '  +1.09%  '
$ws.Range("D49").Value = '''48.92'
$ws.Range("E49").Value = '  -8.34%  '
$ws.Range("D50").Value = '''93.01'
$ws.Range("E50").Value = '  +5.55%  '
$ws.Range("E51").Value = '  +0.34%  '
